$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits on the
#    "Default style theme font" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Turn the last (till-now empty) paragraph into a paragraph that
#    carries a theme-based shading fill plus the text "Paragraph shadow",
#    and move the _GoBack bookmark onto it.
$paragraphCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paragraphCount)
$lastRange = $lastPara.Range

$lastRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00743450' w:rsidRPr='00F47343' w:rsidRDefault='00743450' w:rsidP='005E68D6'><w:pPr><w:shd w:val='clear' w:color='auto' w:fill='9D360E' w:themeFill='text2'/></w:pPr><w:r><w:t>Paragraph shadow</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")
